$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "I have a client who looking to hire a Chinese national with a recent masters in engineering. Is he eligible for support?"
$ws.Range("B3").Value = "Is the grant available to a startup?"

$ws.Range("B3").Select()
